$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.430.91'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.348.14'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.60%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.67'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '184.64'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +0.70%  '
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('E10').Value = '  +0.90%  '
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.932.67'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.81%  '
$ws.Range('E13').Value = '  -0.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.40'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.57%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '67.590.45'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000168'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.369.55'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.71%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '446.20'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.86%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.64'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.68'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.57%  '
$ws.Range('E21').Value = '  +2.89%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '74.03'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.09%  '
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.489.77'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.513'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000121'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.194'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.08'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.53%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.98'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '23.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.37'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.31%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.81'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.24'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('E36').Value = '  +5.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '161.69'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '27.61'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.83%  '
$ws.Range('E39').Value = '  -1.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.833.53'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.794'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.49'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.23'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.35'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '24.65'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.38'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '324.89'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0274'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.991'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '31.15'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.21%  '
